$wb = $excel.ActiveWorkbook
$survey = $wb.Worksheets.Item("survey")
$calc = $wb.Worksheets.Item("calculates")

# Insert a new 3-row "if / note / end if" warning block above the
# existing survey content (new rows 7-9; old rows 7-18 shift to 10-21).
$survey.Range("A7:A9").EntireRow.Insert()

# Insert a new row for the "test" calculation on the calculates sheet
# (new row 6).
$calc.Range("A6:B6").EntireRow.Insert()

# Fill in the new cells (order chosen to match the authoring order so
# new shared strings are appended in the same sequence).
$survey.Cells.Item(8, 7).Value = "{{calculates.test}}"
$calc.Cells.Item(6, 1).Value = "test"
$survey.Cells.Item(7, 3).Value = "data('ADA') == null"
$calc.Cells.Item(6, 2).Value = "freebase.echo('Fill in the date')"

$survey.Cells.Item(7, 2).Value = "if"
$survey.Cells.Item(8, 4).Value = "note"
$survey.Cells.Item(9, 2).Value = "end if"

# Leave the UI selection the way it naturally ends up after typing the
# new calculation row, then come back to the survey sheet/new block.
$calc.Activate()
[void]$calc.Range("B7").Select()

$survey.Activate()
[void]$survey.Range("B8").Select()
